# Auto-generated: refresh FFXIV leve market-price figures (columns H-N)
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets, per scheduled market-data run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6870.5
$ws.Range("J62").Value = 9601.200000000001
$ws.Range("L62").Value = 9601.200000000001
$ws.Range("N62").Value = -10849.2

$ws.Range("H65").Value = 6870.5
$ws.Range("J65").Value = 9601.200000000001
$ws.Range("L65").Value = 48006
$ws.Range("N65").Value = -54246

$ws.Range("H112").Value = 2333.4
$ws.Range("J112").Value = 1389
$ws.Range("L112").Value = 4167
$ws.Range("N112").Value = -6383

$ws.Range("H132").Value = 9320
$ws.Range("I132").Value = 9346.666999999999
$ws.Range("K132").Value = 28040.001
$ws.Range("M132").Value = -25510.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 14972.8125
$ws.Range("J44").Value = 14972.8125
$ws.Range("L44").Value = 14972.8125
$ws.Range("N44").Value = -15948.8125

$ws.Range("H45").Value = 2982.5833
$ws.Range("I45").Value = 1361
$ws.Range("K45").Value = 1361
$ws.Range("M45").Value = -984

$ws.Range("H55").Value = 64666.332
$ws.Range("J55").Value = 74499.5
$ws.Range("L55").Value = 74499.5
$ws.Range("N55").Value = -75129.5

$ws.Range("H88").Value = 1807
$ws.Range("J88").Value = 2258
$ws.Range("L88").Value = 2258
$ws.Range("N88").Value = -3070

$ws.Range("H91").Value = 1807
$ws.Range("J91").Value = 2258
$ws.Range("L91").Value = 2258
$ws.Range("N91").Value = -5066

$ws.Range("H97").Value = 746.9
$ws.Range("I97").Value = 718.7778
$ws.Range("K97").Value = 718.7778
$ws.Range("M97").Value = -222.7778

$ws.Range("H132").Value = 4627.643
$ws.Range("I132").Value = 5488.2
$ws.Range("K132").Value = 16464.6
$ws.Range("M132").Value = -13934.6

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("N133").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6333.3335
$ws.Range("I86").Value = 1666.6666
$ws.Range("J86").Value = 8666.666999999999
$ws.Range("K86").Value = 1666.6666
$ws.Range("L86").Value = 8666.666999999999
$ws.Range("M86").Value = -543.6666
$ws.Range("N86").Value = -10912.667

$ws.Range("H89").Value = 6333.3335
$ws.Range("I89").Value = 1666.6666
$ws.Range("J89").Value = 8666.666999999999
$ws.Range("K89").Value = 8333.333000000001
$ws.Range("L89").Value = 43333.335
$ws.Range("M89").Value = -2717.333000000001
$ws.Range("N89").Value = -54565.335

$ws.Range("H134").Value = 1558.5
$ws.Range("I134").Value = 1558.5
$ws.Range("K134").Value = 4675.5
$ws.Range("M134").Value = -2140.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2180
$ws.Range("I122").Value = 1710
$ws.Range("K122").Value = 5130
$ws.Range("M122").Value = -2680

$ws.Range("H132").Value = 2833
$ws.Range("I132").Value = 2833
$ws.Range("K132").Value = 8499
$ws.Range("M132").Value = -5969

$ws.Range("H134").Value = 1404.125
$ws.Range("I134").Value = 1404.125
$ws.Range("K134").Value = 4212.375
$ws.Range("M134").Value = -1677.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 2810
$ws.Range("J82").Value = 2810
$ws.Range("L82").Value = 8430
$ws.Range("N82").Value = -9242

$ws.Range("H85").Value = 2810
$ws.Range("J85").Value = 2810
$ws.Range("L85").Value = 8430
$ws.Range("N85").Value = -11238

$ws.Range("H131").Value = 1809.3889
$ws.Range("J131").Value = 2873.3333
$ws.Range("L131").Value = 8619.999899999999
$ws.Range("N131").Value = -18699.9999

$ws.Range("H139").Value = 1800.5555
$ws.Range("I139").Value = 1025.625
$ws.Range("K139").Value = 3076.875
$ws.Range("M139").Value = 2063.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 371.3846
$ws.Range("I2").Value = 352.66666
$ws.Range("J2").Value = 387.42856
$ws.Range("K2").Value = 352.66666
$ws.Range("L2").Value = 387.42856
$ws.Range("M2").Value = -239.66666
$ws.Range("N2").Value = -613.4285600000001

$ws.Range("H58").Value = 27500
$ws.Range("J58").Value = 25000
$ws.Range("L58").Value = 25000
$ws.Range("N58").Value = -25554

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0

$ws.Range("H97").Value = 692.7273
$ws.Range("J97").Value = 1025
$ws.Range("L97").Value = 1025
$ws.Range("N97").Value = -2017

$ws.Range("H102").Value = 1939.6471
$ws.Range("I102").Value = 1732.6
$ws.Range("K102").Value = 1732.6
$ws.Range("M102").Value = -110.5999999999999

$ws.Range("H122").Value = 693
$ws.Range("I122").Value = 689
$ws.Range("J122").Value = 695
$ws.Range("K122").Value = 2067
$ws.Range("L122").Value = 2085
$ws.Range("M122").Value = 383
$ws.Range("N122").Value = -6985

$ws.Range("H126").Value = 6269.3335
$ws.Range("I126").Value = 6129.4
$ws.Range("J126").Value = 6969
$ws.Range("K126").Value = 18388.2
$ws.Range("L126").Value = 20907
$ws.Range("M126").Value = -15918.2
$ws.Range("N126").Value = -25847

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7107.8335
$ws.Range("J7").Value = 6500
$ws.Range("L7").Value = 6500
$ws.Range("N7").Value = -6724

$ws.Range("H14").Value = 2005
$ws.Range("J14").Value = 2005
$ws.Range("L14").Value = 2005
$ws.Range("N14").Value = -2349

$ws.Range("H22").Value = 1475.7142
$ws.Range("I22").Value = 1716
$ws.Range("K22").Value = 1716
$ws.Range("M22").Value = -1421

$ws.Range("H27").Value = 1475.7142
$ws.Range("I27").Value = 1716
$ws.Range("K27").Value = 1716
$ws.Range("M27").Value = -1609

$ws.Range("H40").Value = 3642.0715
$ws.Range("I40").Value = 3349.5454
$ws.Range("J40").Value = 4714.6665
$ws.Range("K40").Value = 3349.5454
$ws.Range("L40").Value = 4714.6665
$ws.Range("M40").Value = -3213.5454
$ws.Range("N40").Value = -4986.6665

$ws.Range("H86").Value = 60000
$ws.Range("J86").Value = 60000
$ws.Range("L86").Value = 60000
$ws.Range("N86").Value = -62372

$ws.Range("H89").Value = 60000
$ws.Range("J89").Value = 60000
$ws.Range("L89").Value = 180000
$ws.Range("N89").Value = -191856

$ws.Range("H126").Value = 7107.8335
$ws.Range("J126").Value = 6500
$ws.Range("L126").Value = 19500
$ws.Range("N126").Value = -24440

$ws.Range("H136").Value = 1674.75
$ws.Range("J136").Value = 2499.5
$ws.Range("L136").Value = 7498.5
$ws.Range("N136").Value = -12598.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 7533.3335
$ws.Range("J3").Value = 8800
$ws.Range("L3").Value = 8800
$ws.Range("N3").Value = -9028

$ws.Range("H62").Value = 9889.223
$ws.Range("I62").Value = 5668
$ws.Range("J62").Value = 11999.833
$ws.Range("K62").Value = 5668
$ws.Range("L62").Value = 11999.833
$ws.Range("M62").Value = -5044
$ws.Range("N62").Value = -13247.833

$ws.Range("H65").Value = 9889.223
$ws.Range("I65").Value = 5668
$ws.Range("J65").Value = 11999.833
$ws.Range("K65").Value = 28340
$ws.Range("L65").Value = 59999.165
$ws.Range("M65").Value = -25220
$ws.Range("N65").Value = -66239.16500000001

$ws.Range("H107").Value = 670.4
$ws.Range("I107").Value = 784
$ws.Range("K107").Value = 2352
$ws.Range("M107").Value = -432

$ws.Range("H126").Value = 4457.387
$ws.Range("I126").Value = 2852.0527
$ws.Range("K126").Value = 8556.158100000001
$ws.Range("M126").Value = -6086.158100000001

$ws.Range("H136").Value = 2489.2354
$ws.Range("I136").Value = 1758.7858
$ws.Range("K136").Value = 5276.357400000001
$ws.Range("M136").Value = -2726.357400000001
